# Update participant/attendance-style numbers on the "展览" (sheet 1) and
# "全部类型" (sheet 4) worksheets. Both sheets mirror the same underlying
# rows (sheet 4 aggregates all categories, so its row numbers are shifted
# by the extra rows contributed by other categories), so each logical
# record is updated on both sheets at its own row.

$wb = $excel.ActiveWorkbook

$wsExhibit = $wb.Worksheets.Item(1)   # 展览
$wsAll     = $wb.Worksheets.Item(4)   # 全部类型

# 展览 sheet (rows match directly)
$wsExhibit.Range("G2").Value  = 60
$wsExhibit.Range("F3").Value  = 100
$wsExhibit.Range("F4").Value  = 70
$wsExhibit.Range("F7").Value  = 6998
$wsExhibit.Range("F18").Value = 51
$wsExhibit.Range("F20").Value = 5358
$wsExhibit.Range("F21").Value = 132
$wsExhibit.Range("F22").Value = 190
$wsExhibit.Range("F23").Value = 787
$wsExhibit.Range("F25").Value = 265

# 全部类型 sheet (row numbers shifted after row 18 due to extra entries)
$wsAll.Range("G2").Value  = 60
$wsAll.Range("F3").Value  = 100
$wsAll.Range("F4").Value  = 70
$wsAll.Range("F7").Value  = 6998
$wsAll.Range("F18").Value = 51
$wsAll.Range("F21").Value = 5358
$wsAll.Range("F23").Value = 132
$wsAll.Range("F24").Value = 190
$wsAll.Range("F25").Value = 787
$wsAll.Range("F27").Value = 265
